$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source values in C/D/E for rows 2, 3 and 5 are stored as text
# (numeric-looking strings). Force the number format to Text first so the
# new values keep the same "stored as text" semantics as the rest of the
# sheet instead of being auto-coerced to numbers.
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C3:E3").NumberFormat = "@"
$ws.Range("C5:E5").NumberFormat = "@"

# Row 2: runs, balls, fours
$ws.Range("C2").Value = "11"
$ws.Range("D2").Value = "15"
$ws.Range("E2").Value = "1"

# Row 3: runs, balls, fours
$ws.Range("C3").Value = "24"
$ws.Range("D3").Value = "20"
$ws.Range("E3").Value = "3"

# Row 5: runs, balls, fours
$ws.Range("C5").Value = "5"
$ws.Range("D5").Value = "8"
$ws.Range("E5").Value = "0"
